$d = $word.ActiveDocument

# --- Step 1: merge "cd backend; " and "node " runs in the first backend/node paragraph ---
$null = $d.Content.Find.Execute("cd backend; node ", $true, $false, $false, $false, $false, $true, 1, $false, "cd backend; node ", 2)

# --- Step 2: append the new "db migration tips" section as 8 new paragraphs at the end of the body ---
$targetParas = @(
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@,
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>db</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve"> migration tips</w:t>
      </w:r>
    </w:p>
'@,
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@,
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>cd backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>/app</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">; </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>npx</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>sequelize</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">-cli </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>migration:create</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> --name </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>modify_users_add_new_fields</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
'@,
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@,
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:i/>
          <w:iCs/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>edit</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> .</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>js</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> file under migrations/ &#8230;</w:t>
      </w:r>
    </w:p>
'@,
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@,
@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>cd backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>/app</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">; </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>npx</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>sequelize</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">-cli </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:sz w:val="22"/>
          <w:szCs w:val="22"/>
        </w:rPr>
        <w:t>db:migrate</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
'@
)

foreach ($fragXml in $targetParas) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter()
    $newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newLast.Range.InsertXML($fragXml)
}

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
